# Logged Week 15 and simulated Week 16
$wb = $excel.ActiveWorkbook

# --- Rushing sheet ---
$wsRushing = $wb.Worksheets.Item("Rushing")

# Row 2 (J.Hurts)
$wsRushing.Range("C2").Value = 42
$wsRushing.Range("D2").Value = 43
$wsRushing.Range("E2").Value = 34
$wsRushing.Range("F2").Value = 29

# Row 4 (M.Sanders)
$wsRushing.Range("C4").Value = 42
$wsRushing.Range("D4").Value = 21
$wsRushing.Range("E4").Value = 4
$wsRushing.Range("F4").Value = 9

# Row 7 (J.Howard)
$wsRushing.Range("C7").Value = 31
$wsRushing.Range("D7").Value = 26
$wsRushing.Range("E7").Value = 10

# --- Receiving sheet ---
$wsReceiving = $wb.Worksheets.Item("Receiving")

# Row 2 (M.Sanders)
$wsReceiving.Range("C2").Value = 6
$wsReceiving.Range("D2").Value = 6

# Row 4 (K.Gainwell)
$wsReceiving.Range("C4").Value = 40

# Row 5 (D.Smith)
$wsReceiving.Range("C5").Value = 57
$wsReceiving.Range("D5").Value = 41
$wsReceiving.Range("E5").Value = 29
$wsReceiving.Range("F5").Value = 12
$wsReceiving.Range("G5").Value = 6

# Row 6 (J.Reagor)
$wsReceiving.Range("C6").Value = 37
$wsReceiving.Range("D6").Value = 25
$wsReceiving.Range("E6").Value = 12
$wsReceiving.Range("F6").Value = 4
$wsReceiving.Range("G6").Value = 5
$wsReceiving.Range("H6").Value = 3

# Row 7 (Q.Watkins)
$wsReceiving.Range("C7").Value = 30
$wsReceiving.Range("D7").Value = 23

# Row 8 (G.Ward)
$wsReceiving.Range("C8").Value = 8
$wsReceiving.Range("D8").Value = 4
$wsReceiving.Range("G8").Value = 7
$wsReceiving.Range("H8").Value = 3

# Row 10 (D.Goedert)
$wsReceiving.Range("C10").Value = 47
$wsReceiving.Range("D10").Value = 33
$wsReceiving.Range("E10").Value = 18
$wsReceiving.Range("F10").Value = 15

# Row 11 (J.Stoll)
$wsReceiving.Range("C11").Value = 4
$wsReceiving.Range("D11").Value = 3
